{"js": "// Load the body with its paragraphs and tables so we can inspect/modify them.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nlet tables = body.tables;\nparagraphs.load(\"items,style,text\");\ntables.load(\"items\");\nawait context.sync();\n\n// --- 1. Convert the first paragraph (title) from Heading1 \"NIC Claim\n//        Template\" to Title style \"Claim Report Template\". ---\nconst titlePara = paragraphs.items[0];\ntitlePara.style = \"Title\";\ntitlePara.clear();\ntitlePara.insertText(\"Claim Report Template\", Word.InsertLocation.start);\nawait context.sync();\n\n// --- 2. Replace the field table with plain \"Label: {{placeholder}}\"\n//        paragraphs, in the same order as the table rows. ---\nconst rows = [\n  [\"Claim Number\", \"{{claim_no}}\"],\n  [\"Patient Name\", \"{{patient_name}}\"],\n  [\"Policy Number\", \"{{Policyno}}\"],\n  [\"Date of Admission\", \"{{doa}}\"],\n  [\"Date of Discharge\", \"{{dod}}\"],\n  [\"Insured Name\", \"{{insured_name}}\"],\n  [\"Hospital Name\", \"{{hospital_name}}\"],\n  [\"City\", \"{{city}}\"],\n  [\"State\", \"{{state}}\"],\n];\n\nlet table = tables.items[0];\n\n// Insert the replacement paragraphs right after the table (advancing the\n// anchor each time keeps them in the original top-to-bottom row order).\nlet anchor = table.getRange(Word.RangeLocation.after);\nfor (const [label, placeholder] of rows) {\n  anchor = anchor.insertParagraph(`${label}: ${placeholder}`, Word.InsertLocation.after);\n}\nawait context.sync();\n\n// Re-load the table reference before deleting it; deleting it via a stale\n// anchor (one obtained before the inserts above) is a silent no-op.\ntables = body.tables;\ntables.load(\"items\");\nawait context.sync();\ntable = tables.items[0];\ntable.delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Replace the field table with plain \"Label: {{placeholder}}\"\n#        paragraphs, in the same order as the table rows. Do this while the\n#        title paragraph is untouched so the new paragraphs are created next\n#        to the table (picking up no explicit paragraph style) instead of\n#        inheriting the title's Heading1 style. ---\n$t = $d.Tables.Item(1)\n\n$rows = @(\n    @(\"Claim Number\", \"{{claim_no}}\"),\n    @(\"Patient Name\", \"{{patient_name}}\"),\n    @(\"Policy Number\", \"{{Policyno}}\"),\n    @(\"Date of Admission\", \"{{doa}}\"),\n    @(\"Date of Discharge\", \"{{dod}}\"),\n    @(\"Insured Name\", \"{{insured_name}}\"),\n    @(\"Hospital Name\", \"{{hospital_name}}\"),\n    @(\"City\", \"{{city}}\"),\n    @(\"State\", \"{{state}}\")\n)\n\n$pos = $t.Range.End\nforeach ($row in $rows) {\n    $label = $row[0]\n    $placeholder = $row[1]\n    $text = \"$($label): $($placeholder)\"\n\n    $r = $d.Range($pos, $pos)\n    $r.InsertParagraphAfter()\n\n    $r2 = $d.Range($pos, $pos)\n    $r2.InsertAfter($text)\n\n    $pos = $pos + $text.Length + 1\n}\n\n# Re-fetch the table handle before deleting it \u2014 a handle obtained before the\n# structural edits above can be stale.\n$t2 = $d.Tables.Item(1)\n$t2.Delete()\n\n# --- 2. Convert the first paragraph (title) from Heading1 \"NIC Claim\n#        Template\" to Title style \"Claim Report Template\". Re-fetch via\n#        Content.Paragraphs so we see the post-edit document state. ---\n$titlePara = $d.Content.Paragraphs.Item(1)\n$titlePara.Range.Text = \"Claim Report Template\"\n$titlePara.Style = $d.Styles.Item(\"Title\")\n"}
